# Update parts for availability
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new BOM line for R26 (now RC0201JR-0710KL / 10K), pushing the
#            "R10, R12, R14" group down to row 7 and the big R5.. group to row 8 ---
$ws.Range("D6").Value = "R26"
$ws.Range("E6").Value = "RC0201JR-0710KL"
$ws.Range("F6").Value = "10K"
$ws.Range("H6").Value = 1

# --- Row 7: becomes the former "R10, R12, R14 / RC0603FR-07100KL / 100K" line ---
$ws.Range("D7").Value = "R10, R12, R14"
$ws.Range("E7").Value = "RC0603FR-07100KL"
$ws.Range("F7").Value = "100K"
$ws.Range("H7").Value = 3

# --- Row 8: becomes the former big "R5, R6, ... / RT0603FRE071KL / 1K" line
#            (replaces the old R26 / SR0603KR-7W10KL / 10K line, now unavailable) ---
$ws.Range("D8").Value = "R5, R6, R7, R8, R9, R11, R13, R19, R20, R21, R22, R23, R24, R25"
$ws.Range("E8").Value = "RT0603FRE071KL"
$ws.Range("F8").Value = "1K"
$ws.Range("H8").Value = 14

# --- Row 13 (Q1 crystal): swap the manf# for an available alternative ---
$ws.Range("E13").Value = "ECS-160-10-36-CKM-TR"

# --- Row 15 (CT.. current-transformer jacks): swap the manf# for an available alternative ---
$ws.Range("E15").Value = "ST-PJ-342"

# --- Widen the sheet's columns (so the longer reference/manf# text fits) ---
$ws.Range("A1:AMK1").EntireColumn.ColumnWidth = 29.666666666666668

# --- Update the saved cursor / scroll position ---
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H11").Select() | Out-Null
